$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.071777105331421
$ws.Range("B1").Value = 2.602810621261597
$ws.Range("C1").Value = 3.698626041412354
$ws.Range("D1").Value = 5.696628093719482
$ws.Range("E1").Value = 1.71710991859436
